# Updated cryptos list (prices + 1h volume %) to match the latest scrape.
# Values are stored as plain text in the source sheet (Price column mixes
# thousand-separator-as-dot numbers like "27.679.30" with plain decimals
# like "207.16", and Volume is a padded "  +x.xx%  " string), so for any
# new value that Excel would otherwise auto-detect as a real number we
# force the cell to Text first and clear the number format back off again
# afterwards so no stray style index is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.679.30"
$ws.Range("D3").Value = "1.588.43"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  -4.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.253"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0869"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "1.814.19"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").Value = "1.612.50"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.86"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("D16").Value = "27.664.33"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.43"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.37"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.06%  "
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.97%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("E22").Value = "  -5.05%  "
$ws.Range("E23").Value = "  -3.62%  "
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.46"
$ws.Range("D25").ClearFormats()
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.13"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("E29").Value = "  -4.60%  "
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("E32").Value = "  -5.09%  "
$ws.Range("D33").Value = "1.370.19"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.79%  "
$ws.Range("E35").Value = "  -4.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.981"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.973"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("E45").Value = "  -3.73%  "
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("D47").Value = "1.725.47"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0967"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.33%  "
